# Applies the "456a3b4" content refresh to the 合肥-漫展信息 workbook:
#   - bumps several "want to go" counters (F column) on both the
#     展览 (sheet 1) and 全部类型 (sheet 4) sheets
#   - inserts a brand-new exhibition row ("次元之门...") as the new
#     row 15 on both sheets, pushing every later row down by one
#   - bumps two more F-column counters that land on the rows that were
#     shifted down

$wb = $excel.ActiveWorkbook

function Update-MangaSheet($ws, $lastOldDataRow) {
    # $lastOldDataRow = last populated data row BEFORE the insert

    # ---- 1. shift rows 15..lastOldDataRow down by one (bottom-up so we
    #         never clobber a row before it has been read). Column A is
    #         a plain running index (row# - 1), so it is rewritten from
    #         scratch afterwards instead of being shifted along with the
    #         rest of the row. ---------------------------------------
    for ($r = $lastOldDataRow; $r -ge 15; $r--) {
        $src = $ws.Range("B" + $r + ":I" + $r)
        $dst = $ws.Range("B" + ($r + 1) + ":I" + ($r + 1))
        $src.Copy()
        $dst.PasteSpecial(-4104)   # xlPasteAll
    }
    $excel.CutCopyMode = 0

    # the very last destination row sits outside the sheet's original
    # used range, so the paste above leaves column A unstyled there;
    # borrow the style from a still-intact index cell (column A is
    # uniformly styled for every data row) before filling in values.
    $ws.Range("A14").Copy()
    $ws.Range("A" + ($lastOldDataRow + 1)).PasteSpecial(-4122)   # xlPasteFormats
    $excel.CutCopyMode = 0

    for ($r = 15; $r -le ($lastOldDataRow + 1); $r++) {
        $ws.Range("A" + $r).Value2 = ($r - 1)
    }

    # ---- 2. write the brand-new row 15 -------------------------------
    $ws.Range("B15").Value2 = "2024-05-01"
    $ws.Range("C15").Value2 = "合肥·第十三届次元之门动漫游戏博览会·触手猫X福瑞福瑞福兽漫联合专区"
    $ws.Range("D15").Value2 = "南京路与庐州大道交汇处 合肥滨湖国际会展中心"
    $ws.Range("E15").Value2 = "2024.05.01 10:00-05.03 17:00"
    $ws.Range("F15").Value2 = 0
    $ws.Range("G15").Value2 = 89
    $ws.Range("H15").Value2 = "https://show.bilibili.com/platform/detail.html?id=83770"
    $ws.Range("I15").Value2 = "//i2.hdslb.com/bfs/openplatform/202404/91GGPKXa1712065132816.jpeg"

    # ---- 3. refresh the "want to go" counters ------------------------
    $ws.Range("F2").Value2 = 8854
    $ws.Range("F3").Value2 = 8266
    $ws.Range("F7").Value2 = 47
    $ws.Range("F10").Value2 = 207
    $ws.Range("F12").Value2 = 749
    $ws.Range("F14").Value2 = 5336

    # rows that used to be 20/21 (运动番only-群青日和 / 第六届环形宇宙…)
    # are now 21/22 after the insert above; their counters also moved up.
    $ws.Range("F" + 21).Value2 = 156
    $ws.Range("F" + 22).Value2 = 157
}

# 展览 sheet: originally had data rows 2..22 (row 22 = 安徽·MAX特摄only展)
$wsExhibit = $wb.Worksheets.Item("展览")
Update-MangaSheet $wsExhibit 22

# 全部类型 sheet: identical layout, but one extra trailing row
# (合肥·首届包河留声机音乐节…), so the last old data row is 23
$wsAll = $wb.Worksheets.Item("全部类型")
Update-MangaSheet $wsAll 23
